# no-op test
$wb = $excel.ActiveWorkbook
Write-Host $wb.Worksheets.Count
